$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $addr, $val)
    $esc = $val.Replace('"', '""')
    $ws.Range("ZZ1").Formula = '="' + $esc + '"'
    $ws.Range("ZZ1").Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4163) | Out-Null
    $ws.Range("ZZ1").Clear() | Out-Null
}

Set-TextValue $ws "D2" '37.821.32'
Set-TextValue $ws "E2" '  -0.83%  '
Set-TextValue $ws "D3" '2.035.70'
Set-TextValue $ws "E4" '  +0.04%  '
Set-TextValue $ws "D5" '228.22'
Set-TextValue $ws "E5" '  -0.95%  '
Set-TextValue $ws "E6" '  -1.37%  '
Set-TextValue $ws "D7" '60.30'
Set-TextValue $ws "E7" '  +1.50%  '
Set-TextValue $ws "E8" '  -0.06%  '
Set-TextValue $ws "E9" '  -2.00%  '
Set-TextValue $ws "D10" '0.0822'
Set-TextValue $ws "E10" '  +1.10%  '
Set-TextValue $ws "D11" '0.104'
Set-TextValue $ws "E11" '  -0.40%  '
Set-TextValue $ws "D12" '2.337.35'
Set-TextValue $ws "E12" '  -1.18%  '
Set-TextValue $ws "D13" '14.49'
Set-TextValue $ws "E13" '  -1.92%  '
Set-TextValue $ws "D14" '21.09'
Set-TextValue $ws "E14" '  -0.89%  '
Set-TextValue $ws "E15" '  +0.37%  '
Set-TextValue $ws "D16" '5.19'
Set-TextValue $ws "E16" '  -2.12%  '
Set-TextValue $ws "D17" '2.070.75'
Set-TextValue $ws "E17" '  -0.22%  '
Set-TextValue $ws "D18" '37.773.15'
Set-TextValue $ws "E18" '  -0.79%  '
Set-TextValue $ws "E19" '  -0.24%  '
Set-TextValue $ws "D20" '5.91'
Set-TextValue $ws "E20" '  -5.79%  '
Set-TextValue $ws "E21" '  -1.45%  '
Set-TextValue $ws "D22" '223.84'
Set-TextValue $ws "E22" '  -0.65%  '
Set-TextValue $ws "E23" '  -0.02%  '
Set-TextValue $ws "D24" '2.41'
Set-TextValue $ws "E24" '  -1.19%  '
Set-TextValue $ws "E25" '  +0.54%  '
Set-TextValue $ws "D26" '9.38'
Set-TextValue $ws "E26" '  +1.04%  '
Set-TextValue $ws "D27" '167.48'
Set-TextValue $ws "E27" '  +0.55%  '
Set-TextValue $ws "E28" '  -2.74%  '
Set-TextValue $ws "E29" '  -0.97%  '
Set-TextValue $ws "D30" '1.27'
Set-TextValue $ws "E30" '  -4.10%  '
Set-TextValue $ws "E31" '  +1.12%  '
Set-TextValue $ws "D32" '2.18'
Set-TextValue $ws "E32" '  +6.02%  '
Set-TextValue $ws "D33" '4.41'
Set-TextValue $ws "E33" '  -3.45%  '
Set-TextValue $ws "E34" '  -0.28%  '
Set-TextValue $ws "D36" '6.32'
Set-TextValue $ws "E36" '  +2.45%  '
Set-TextValue $ws "E37" '  -2.46%  '
Set-TextValue $ws "D38" '3.33'
Set-TextValue $ws "E38" '  +1.30%  '
Set-TextValue $ws "E39" '  +0.13%  '
Set-TextValue $ws "D40" '17.74'
Set-TextValue $ws "E40" '  +4.01%  '
Set-TextValue $ws "D41" '1.536.73'
Set-TextValue $ws "E41" '  +0.09%  '
Set-TextValue $ws "E42" '  -0.22%  '
Set-TextValue $ws "D43" '96.22'
Set-TextValue $ws "E43" '  -2.08%  '
Set-TextValue $ws "E44" '  -2.54%  '
Set-TextValue $ws "E45" '  -1.30%  '
Set-TextValue $ws "E46" '  -2.53%  '
Set-TextValue $ws "D47" '4.02'
Set-TextValue $ws "E47" '  -3.84%  '
Set-TextValue $ws "E48" '  -1.50%  '
Set-TextValue $ws "E49" '  +0.44%  '
Set-TextValue $ws "E50" '  -0.49%  '
Set-TextValue $ws "D51" '2.226.98'
Set-TextValue $ws "E51" '  -1.12%  '
